$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the existing two data rows (2-3) twice more into rows 4-7
$ws.Range("A2:F3").Copy()
$ws.Range("A4").PasteSpecial(-4104)
$ws.Range("A6").PasteSpecial(-4104)
$excel.CutCopyMode = 0

# Update selection to match target state
$ws.Range("A2:G4").Select()
